{"js": "// \"Laravel Inventory Project Step-21 (Live Class 2) :\" is followed by a\n// value \"1.10.00\" which must become \"1.50.00 STOPED\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.text.indexOf(\"Laravel Inventory Project Step-21 (Live Class 2) :\") !== -1\n);\n\nif (target) {\n  const results = target.search(\"1.10.00\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"1.50.00 STOPED\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Find the paragraph \"Laravel Inventory Project Step-21 (Live Class 2) :\"\n# whose value \"1.10.00\" needs to become \"1.50.00 STOPED\".\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Laravel Inventory Project Step-21 (Live Class 2) :*\") {\n        $target = $p\n    }\n}\n\nif ($target -ne $null) {\n    $r = $target.Range\n    $found = $r.Find.Execute(\"1.10.00\")\n    if ($found) {\n        $r.Text = \"1.50.00 STOPED\"\n    }\n}\n"}
